$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set the Runmode column (C) to "Y" for rows 2-12 (already "Y" in row 13)
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 3).Value = "Y"
}

# Select C2:C13 with active cell C13, matching the final selection state
$ws.Range("C2:C13").Select()
$ws.Range("C13").Activate()
